$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the activity category for row 63 (was "Code ", now "Code/Debug ")
$ws.Range("F63").Value = "Code/Debug "

# Add the Stop time for row 63 (typed as text, keeps the h:mm stop-time style)
$ws.Range("C63").NumberFormat = "h:mm"
$ws.Range("C63").Value = "4:18PM"

# Add the follow-up note for row 63
$ws.Range("H63").Value = "Now that adding orders to db is complete, now to add invoices"

# New log entry row 64
$ws.Range("B64").Value = "4:18PM"
$ws.Range("F64").Value = "Code"
$ws.Range("G64").Value = "Generate Invoices GUI page and functionality start"

# Update sheet view to reflect the new active cell / scroll position
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G64").Select()
